$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be auto-detected
# as numbers by Excel, so they stay stored as text (matching the original inlineStr type).
$ws.Range("D2").Value = '67.281.91'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = '3.871.69'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '469.63'
$ws.Range("E5").Value = '  +10.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.57'
$ws.Range("E6").Value = '  +10.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.635'
$ws.Range("E7").Value = '  +3.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000314'
$ws.Range("E11").Value = '  -7.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.52'
$ws.Range("E12").Value = '  +3.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.46'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '4.498.12'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.84'
$ws.Range("E15").Value = '  -6.41%  '
$ws.Range("D16").Value = '3.903.99'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.12'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  +5.98%  '
$ws.Range("D20").Value = '67.447.63'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '437.52'
$ws.Range("E21").Value = '  +5.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.93'
$ws.Range("E22").Value = '  -1.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.37'
$ws.Range("E23").Value = '  +5.15%  '
$ws.Range("E24").Value = '  +5.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.61'
$ws.Range("E25").Value = '  +10.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '37.96'
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.24'
$ws.Range("E27").Value = '  +10.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.49'
$ws.Range("E29").Value = '  +2.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '729.27'
$ws.Range("E30").Value = '  +1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.88'
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  +7.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.77'
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '44.32'
$ws.Range("E34").Value = '  +13.19%  '
$ws.Range("E35").Value = '  +7.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.95'
$ws.Range("E36").Value = '  +3.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.54'
$ws.Range("E38").Value = '  -5.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0483'
$ws.Range("E39").Value = '  +3.88%  '
$ws.Range("E40").Value = '  +9.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '0.0₃0693'
$ws.Range("E43").Value = '  +3.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  +13.08%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.47'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.29'
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  +6.72%  '
$ws.Range("E49").Value = '  +5.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.40'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("E51").Value = '  +1.14%  '
